$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $true, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "2024-04-11 Thursday" "2024-04-12 Friday"

Replace-Text "37×47=" "53×43="
Replace-Text "24×79=" "34×18="
Replace-Text "21×53=" "20×29="
Replace-Text "92×23=" "84×74="
Replace-Text "88×88=" "20×11="

Replace-Text "92×17=" "80×14="
Replace-Text "70×82=" "27×19="
Replace-Text "73×57=" "30×24="
Replace-Text "35×20=" "44×96="
Replace-Text "77×57=" "56×36="

Replace-Text "14×12=" "92×88="
Replace-Text "18×70=" "84×69="
Replace-Text "31×78=" "85×87="
Replace-Text "61×47=" "36×43="
Replace-Text "83×74=" "90×81="

Replace-Text "87×18=" "55×78="
Replace-Text "42×37=" "18×31="
Replace-Text "96×38=" "53×46="
Replace-Text "24×13=" "87×83="
Replace-Text "23×68=" "67×42="

Replace-Text "52×69=" "13×59="
Replace-Text "73×30=" "49×79="
Replace-Text "54×51=" "58×38="
Replace-Text "77×47=" "99×61="
Replace-Text "38×54=" "30×18="
